$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-11 16:44:54"
$wsZhCn.Range("H2").Value = "2016-03-11 16:45:14"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-11 16:44:59"
$wsDeDe.Range("H2").Value = "2016-03-11 16:45:19"
